# Both the "展览" (sheet 1) and "全部类型" (sheet 4) sheets have the same
# two event rows; the dates in column B for rows 2 and 3 switch from
# dash-separated to dot-separated (2024-02-14 -> 2024.02.14,
# 2024-02-18 -> 2024.02.18). We force the cells to Text format before
# writing so Excel doesn't reinterpret the dotted date string as a real
# date value, then clear the temporary formatting so the cell keeps its
# original (default/general) style.

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    $datesRange = $ws.Range("B2:B3")
    $datesRange.NumberFormat = "@"

    $ws.Range("B2").Value = "2024.02.14"
    $ws.Range("B3").Value = "2024.02.18"

    $datesRange.ClearFormats()
}
